$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 54
$ws.Range("I2").Value = 158
$ws.Range("J2").Value = 655
$ws.Range("K2").Value = 4
$ws.Range("L2").Value = 148
$ws.Range("M2").Value = 7
$ws.Range("N2").Value = 90
$ws.Range("P2").Value = 3
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 8
$ws.Range("S2").Value = 60
$ws.Range("T2").Value = 118
$ws.Range("U2").Value = 6
$ws.Range("V2").Value = 950
$ws.Range("X2").Value = 1031
$ws.Range("Y2").Value = 3
$ws.Range("Z2").Value = 19
$ws.Range("AA2").Value = 2
